$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Columns: D = Price (text), E = Volume(1h) (text, padded "  +x.xx%  ").
# Cells whose new value parses as a plain number are forced back to Text
# format first so Excel keeps storing them as literal strings (matching
# the original inlineStr cells) instead of auto-converting to numerics.

$ws.Range("D2").Value = "28.516.46"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.820.33"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.02"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5103"
$ws.Range("E7").Value = "  -6.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3950"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08316"
$ws.Range("E9").Value = "  +8.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.63"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.320"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.542"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "1.815.83"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001140"
$ws.Range("E17").Value = "  +5.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.61"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06654"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.79"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.094"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "28.548.67"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.48"
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.265"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.26"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "2.026.62"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.410"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1097"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.789"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.654"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07082"
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02345"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.219"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.829"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6315"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.27"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.180"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.400"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5936"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.734"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.10"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.988"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.186"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06894"
$ws.Range("E51").Value = "  +0.11%  "
